$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the data table (A1:D42) with updated query results ---
$ws.Cells.Item(1,1).Value = "Transformation"
$ws.Cells.Item(1,2).Value = "Method"
$ws.Cells.Item(1,3).Value = "AUC"
$ws.Cells.Item(1,4).Value = "Accuracy"
$ws.Cells.Item(2,1).Value = "filter_and_pca_high-pass"
$ws.Cells.Item(2,2).Value = "log_lasso"
$ws.Cells.Item(2,3).Value = 0.94001095981026195
$ws.Cells.Item(2,4).Value = 0.84907571537098003
$ws.Cells.Item(3,1).Value = "filter_and_pca_high-pass"
$ws.Cells.Item(3,2).Value = "neuralNet"
$ws.Cells.Item(3,3).Value = 0.93789629163983201
$ws.Cells.Item(3,4).Value = 0.83742719675867305
$ws.Cells.Item(4,1).Value = "pca"
$ws.Cells.Item(4,2).Value = "xgboost"
$ws.Cells.Item(4,3).Value = 0.93403288148510999
$ws.Cells.Item(4,4).Value = 0.85110154469485899
$ws.Cells.Item(5,1).Value = "pca"
$ws.Cells.Item(5,2).Value = "neuralNet"
$ws.Cells.Item(5,3).Value = 0.93284882138385306
$ws.Cells.Item(5,4).Value = 0.84704988604709996
$ws.Cells.Item(6,1).Value = "pca"
$ws.Cells.Item(6,2).Value = "log_lasso"
$ws.Cells.Item(6,3).Value = 0.92544838155341602
$ws.Cells.Item(6,4).Value = 0.85363383134970805
$ws.Cells.Item(7,1).Value = "filter_and_pca_low-pass"
$ws.Cells.Item(7,2).Value = "log_lasso"
$ws.Cells.Item(7,3).Value = 0.91232331542991296
$ws.Cells.Item(7,4).Value = 0.83869334008609697
$ws.Cells.Item(8,1).Value = "filter_highpass"
$ws.Cells.Item(8,2).Value = "log_lasso"
$ws.Cells.Item(8,3).Value = 0.91043077087422697
$ws.Cells.Item(8,4).Value = 0.83844011142061203
$ws.Cells.Item(9,1).Value = "filter_and_pca_low-pass"
$ws.Cells.Item(9,2).Value = "neuralNet"
$ws.Cells.Item(9,3).Value = 0.89051604066325796
$ws.Cells.Item(9,4).Value = 0.81337047353760406
$ws.Cells.Item(10,1).Value = "nested"
$ws.Cells.Item(10,2).Value = "log_lasso"
$ws.Cells.Item(10,3).Value = 0.86355806593758699
$ws.Cells.Item(10,4).Value = 0.64168143833881996
$ws.Cells.Item(11,1).Value = "nested"
$ws.Cells.Item(11,2).Value = "randomForest"
$ws.Cells.Item(11,3).Value = 0.84764438548539001
$ws.Cells.Item(11,4).Value = 0.72651304127627203
$ws.Cells.Item(12,1).Value = "filter_highpass"
$ws.Cells.Item(12,2).Value = "xgboost"
$ws.Cells.Item(12,3).Value = 0.84404071872525099
$ws.Cells.Item(12,4).Value = 0.73841478855406395
$ws.Cells.Item(13,1).Value = "gs-ref_highpass"
$ws.Cells.Item(13,2).Value = "log_lasso"
$ws.Cells.Item(13,3).Value = 0.83970712552012805
$ws.Cells.Item(13,4).Value = 0.77184097239807503
$ws.Cells.Item(14,1).Value = "base_features"
$ws.Cells.Item(14,2).Value = "xgboost"
$ws.Cells.Item(14,3).Value = 0.837973148978436
$ws.Cells.Item(14,4).Value = 0.73942770321600404
$ws.Cells.Item(15,1).Value = "filter_lowpass"
$ws.Cells.Item(15,2).Value = "log_lasso"
$ws.Cells.Item(15,3).Value = 0.83765151912009195
$ws.Cells.Item(15,4).Value = 0.800202582932388
$ws.Cells.Item(16,1).Value = "nested"
$ws.Cells.Item(16,2).Value = "xgboost"
$ws.Cells.Item(16,3).Value = 0.83744082267396702
$ws.Cells.Item(16,4).Value = 0.72068878197011899
$ws.Cells.Item(17,1).Value = "filter_lowpass"
$ws.Cells.Item(17,2).Value = "xgboost"
$ws.Cells.Item(17,3).Value = 0.82753783291579197
$ws.Cells.Item(17,4).Value = 0.68473031147125796
$ws.Cells.Item(18,1).Value = "gs-ref_highpass"
$ws.Cells.Item(18,2).Value = "xgboost"
$ws.Cells.Item(18,3).Value = 0.82537386100659704
$ws.Cells.Item(18,4).Value = 0.72853887060015199
$ws.Cells.Item(19,1).Value = "nested"
$ws.Cells.Item(19,2).Value = "neuralNet"
$ws.Cells.Item(19,3).Value = 0.82213971549688403
$ws.Cells.Item(19,4).Value = 0.68473031147125796
$ws.Cells.Item(20,1).Value = "gbf_highpass"
$ws.Cells.Item(20,2).Value = "xgboost"
$ws.Cells.Item(20,3).Value = 0.82147372983810896
$ws.Cells.Item(20,4).Value = 0.75284882248670504
$ws.Cells.Item(21,1).Value = "filter_lowpass"
$ws.Cells.Item(21,2).Value = "randomForest"
$ws.Cells.Item(21,3).Value = 0.80702760598506296
$ws.Cells.Item(21,4).Value = 0.72625981261078698
$ws.Cells.Item(22,1).Value = "gbf_highpass"
$ws.Cells.Item(22,2).Value = "randomForest"
$ws.Cells.Item(22,3).Value = 0.80219776551347199
$ws.Cells.Item(22,4).Value = 0.749303621169916
$ws.Cells.Item(23,1).Value = "pca"
$ws.Cells.Item(23,2).Value = "randomForest"
$ws.Cells.Item(23,3).Value = 0.79813174780726703
$ws.Cells.Item(23,4).Value = 0.70194986072423304
$ws.Cells.Item(24,1).Value = "gs-ref_lowpass"
$ws.Cells.Item(24,2).Value = "xgboost"
$ws.Cells.Item(24,3).Value = 0.79717956935238499
$ws.Cells.Item(24,4).Value = 0.68979488478095696
$ws.Cells.Item(25,1).Value = "gs-ref_lowpass"
$ws.Cells.Item(25,2).Value = "log_lasso"
$ws.Cells.Item(25,3).Value = 0.79498580976768696
$ws.Cells.Item(25,4).Value = 0.664218789566979
$ws.Cells.Item(26,1).Value = "filter_and_pca_low-pass"
$ws.Cells.Item(26,2).Value = "xgboost"
$ws.Cells.Item(26,3).Value = 0.79210847438824805
$ws.Cells.Item(26,4).Value = 0.71891618131172397
$ws.Cells.Item(27,1).Value = "gs-ref_lowpass"
$ws.Cells.Item(27,2).Value = "neuralNet"
$ws.Cells.Item(27,3).Value = 0.79175590129802298
$ws.Cells.Item(27,4).Value = 0.664218789566979
$ws.Cells.Item(28,1).Value = "filter_highpass"
$ws.Cells.Item(28,2).Value = "randomForest"
$ws.Cells.Item(28,3).Value = 0.76266322875804904
$ws.Cells.Item(28,4).Value = 0.66523170422891797
$ws.Cells.Item(29,1).Value = "gs-ref_highpass"
$ws.Cells.Item(29,2).Value = "randomForest"
$ws.Cells.Item(29,3).Value = 0.75898252490609097
$ws.Cells.Item(29,4).Value = 0.68219802481640901
$ws.Cells.Item(30,1).Value = "filter_and_pca_high-pass"
$ws.Cells.Item(30,2).Value = "xgboost"
$ws.Cells.Item(30,3).Value = 0.75642925889290502
$ws.Cells.Item(30,4).Value = 0.65383641428209605
$ws.Cells.Item(31,1).Value = "gbf_highpass"
$ws.Cells.Item(31,2).Value = "log_lasso"
$ws.Cells.Item(31,3).Value = 0.73055558123458597
$ws.Cells.Item(31,4).Value = 0.655862243605976
$ws.Cells.Item(32,1).Value = "filter_and_pca_high-pass"
$ws.Cells.Item(32,2).Value = "randomForest"
$ws.Cells.Item(32,3).Value = 0.71675464302552305
$ws.Cells.Item(32,4).Value = 0.50873638895923001
$ws.Cells.Item(33,1).Value = "base_features"
$ws.Cells.Item(33,2).Value = "randomForest"
$ws.Cells.Item(33,3).Value = 0.70942045509405105
$ws.Cells.Item(33,4).Value = 0.58470498860470999
$ws.Cells.Item(34,1).Value = "gbf_lowpass"
$ws.Cells.Item(34,2).Value = "xgboost"
$ws.Cells.Item(34,3).Value = 0.70780909592350505
$ws.Cells.Item(34,4).Value = 0.62243605976196503
$ws.Cells.Item(35,1).Value = "base_features"
$ws.Cells.Item(35,2).Value = "log_lasso"
$ws.Cells.Item(35,3).Value = 0.68683754520793305
$ws.Cells.Item(35,4).Value = 0.66396556090149395
$ws.Cells.Item(36,1).Value = "filter_and_pca_low-pass"
$ws.Cells.Item(36,2).Value = "randomForest"
$ws.Cells.Item(36,3).Value = 0.64401121351908497
$ws.Cells.Item(36,4).Value = 0.54216257280324098
$ws.Cells.Item(37,1).Value = "gbf_lowpass"
$ws.Cells.Item(37,2).Value = "log_lasso"
$ws.Cells.Item(37,3).Value = 0.63399857430022
$ws.Cells.Item(37,4).Value = 0.59913902253735096
$ws.Cells.Item(38,1).Value = "gbf_lowpass"
$ws.Cells.Item(38,2).Value = "randomForest"
$ws.Cells.Item(38,3).Value = 0.59195532464971201
$ws.Cells.Item(38,4).Value = 0.57432261331982704
$ws.Cells.Item(39,1).Value = "gbf_highpass"
$ws.Cells.Item(39,2).Value = "neuralNet"
$ws.Cells.Item(39,3).Value = 0.58745276342319896
$ws.Cells.Item(39,4).Value = 0.59052924791086303
$ws.Cells.Item(40,1).Value = "gs-ref_lowpass"
$ws.Cells.Item(40,2).Value = "randomForest"
$ws.Cells.Item(40,3).Value = 0.57963760724847102
$ws.Cells.Item(40,4).Value = 0.55203849075715306
$ws.Cells.Item(41,1).Value = "gbf_lowpass"
$ws.Cells.Item(41,2).Value = "neuralNet"
$ws.Cells.Item(41,3).Value = 0.50195982361587399
$ws.Cells.Item(41,4).Value = 0.48316029374525199
$ws.Cells.Item(42,1).Value = "gs-ref_highpass"
$ws.Cells.Item(42,2).Value = "neuralNet"
$ws.Cells.Item(42,3).Value = 0.5
$ws.Cells.Item(42,4).Value = 0.48316029374525199
Write-Host "done writing data"
